# Actualización automática 2025-10-06 16:30:09
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# VENTAS POR GRUPO: registrar venta de INODOROS (17.99) para
# BARROS YUNGA DIEGO VINICIO y actualizar el resumen de la fila 14.
$wsGrupo.Range("H3").Value = 17.99
$wsGrupo.Range("H14").Value = "1 de 12"

# VENTA MENSUAL: mover el presupuesto de ASES GAVILANEZ FAUSTO HERNAN a 0,
# registrar la venta de octubre de BARROS YUNGA DIEGO VINICIO,
# y mover el presupuesto de FABIMP BENIGNO BRAVO S.A.S. a 0.
$wsMensual.Range("G2").Value = 0
$wsMensual.Range("F3").Value = 17.99
$wsMensual.Range("G9").Value = 0

# Totales de la fila 14: octubre sube a 17.99, presupuesto baja a 0.
$wsMensual.Range("F14").Value = 17.99
$wsMensual.Range("G14").Value = 0
